# Applies the "faltan Test y HTML" edit to the recuento workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# F10 previously read "FALTA"; it now reads "OK".
$ws.Range("F10").Value = "OK"

# F25 previously read "Faltan todos"; that cell is now cleared entirely.
$ws.Range("F25").Clear()

# Update the on-screen selection to match the saved view (F3:F24, active cell F3).
$ws.Activate()
$ws.Range("F3:F24").Select()
